$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.685899999999996
$ws.Range("A8").Value = -21.0817
$ws.Range("A10").Value = -20.43379999999997
$ws.Range("A12").Value = -22.50090000000003
$ws.Range("B13").Value = 6.480099999999997
$ws.Range("A18").Value = -22.46180000000003
$ws.Range("D20").Value = -8.250900000000001
